$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the player roster table (A2:C19) with the updated roster:
# re-ordered rows, plus two player swaps (Brandon Clarke -> Cole Anthony,
# Dillon Brooks -> Khris Middleton), each carrying their own Position/Team.

$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Cole Anthony", "PG", "Orlando Magic"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
